# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy the "last row" border formatting (currently on row 43) onto the
#     row that will become the new last row of the table (row 36) before we
#     delete the trailing rows. ---
$ws.Range("B43:J43").Copy()
$ws.Range("B36:J36").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- VALOR MORA ---
$ws.Range("E11").Value = 1100579

# --- Cant. Trabajadores / Cant. Periodos ---
$ws.Range("C13").Value = 6
$ws.Range("F13").Value = 11

# --- Replace the worker detail table (B16:G36) with the new dataset ---
$data = @(
    @("CC", "1065875439", "JAN ESLEIDER RINALDY QUINTERO", "2412", 52000, 1300000),
    @("CC", "92131117", "SAMID ANTONIO RUIZ MERIÝO", "2412", 52000, 1300000),
    @("CC", "92131117", "SAMID ANTONIO RUIZ MERIÝO", "2411", 43333, 1300000),
    @("CC", "1046430970", "CLODOMIRO JOSE HERAZO EPALZA", "2507", 56940, 1300000),
    @("CC", "1046430970", "CLODOMIRO JOSE HERAZO EPALZA", "2506", 56940, 1300000),
    @("CC", "1046430970", "CLODOMIRO JOSE HERAZO EPALZA", "2505", 56940, 1300000),
    @("CC", "1046430970", "CLODOMIRO JOSE HERAZO EPALZA", "2504", 56940, 1300000),
    @("CC", "1046430970", "CLODOMIRO JOSE HERAZO EPALZA", "2503", 56940, 1300000),
    @("CC", "1007857666", "LEONARDO JOSE MENDEZ AISLANT", "2507", 52000, 1300000),
    @("CC", "1007857666", "LEONARDO JOSE MENDEZ AISLANT", "2506", 52000, 1300000),
    @("CC", "1007857666", "LEONARDO JOSE MENDEZ AISLANT", "2505", 52000, 1300000),
    @("CC", "1007857666", "LEONARDO JOSE MENDEZ AISLANT", "2504", 52000, 1300000),
    @("CC", "1007857666", "LEONARDO JOSE MENDEZ AISLANT", "2503", 52000, 1300000),
    @("CC", "1007857666", "LEONARDO JOSE MENDEZ AISLANT", "2502", 52000, 1300000),
    @("CC", "1007857666", "LEONARDO JOSE MENDEZ AISLANT", "2501", 52000, 1300000),
    @("CC", "1007857666", "LEONARDO JOSE MENDEZ AISLANT", "2412", 52000, 1300000),
    @("CC", "1007857666", "LEONARDO JOSE MENDEZ AISLANT", "2411", 52000, 1300000),
    @("CC", "1007857666", "LEONARDO JOSE MENDEZ AISLANT", "2410", 52000, 1300000),
    @("CC", "1007857666", "LEONARDO JOSE MENDEZ AISLANT", "2409", 34666, 1300000),
    @("CC", "92131158", "JEAN CARLOS MARTINEZ GASPAR", "2504", 56940, 1300000),
    @("CC", "1098633317", "CRISTIAN YAIR PEREZ MONTENEGRO", "2504", 56940, 1423500)
)

$r = 16
foreach ($row in $data) {
    $ws.Range("B$r").Value = $row[0]
    $ws.Range("C$r").Value = $row[1]
    $ws.Range("D$r").Value = $row[2]
    $ws.Range("E$r").Value = $row[3]
    $ws.Range("F$r").Value = $row[4]
    $ws.Range("G$r").Value = $row[5]
    $r = $r + 1
}

# --- Remove the now-unused trailing table rows (old rows 37-43); this
#     shifts the gap + footer rows (44-49) up to their new positions
#     (37-42) automatically. ---
$ws.Range("A37:A43").EntireRow.Delete()

# --- Column D ("Nombre Trabajador") re-fit to the new (shorter) content ---
$ws.Columns("D").ColumnWidth = 33.43
